# Apply cryptocurrency price/volume updates per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.734.80"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "2.775.83"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'356.22"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").Value = "'108.79"
$ws.Range("E6").Value = "  -1.92%  "
$ws.Range("D7").Value = "'0.555"
$ws.Range("E7").Value = "  -1.34%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D9").Value = "'0.585"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("D10").Value = "'39.65"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("D12").Value = "'0.0843"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "'19.44"
$ws.Range("E13").Value = "  -1.42%  "
$ws.Range("D14").Value = "'7.60"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").Value = "3.211.71"
$ws.Range("E15").Value = "  -1.13%  "
$ws.Range("D16").Value = "2.774.48"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").Value = "'0.932"
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").Value = "51.665.38"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "'7.46"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").Value = "'3.09"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "'13.11"
$ws.Range("E21").Value = "  -1.19%  "
$ws.Range("D22").Value = "0.0₃0968"
$ws.Range("E22").Value = "  -1.89%  "
$ws.Range("D23").Value = "'70.11"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'268.21"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").Value = "'2.73"
$ws.Range("E25").Value = "  -2.75%  "
$ws.Range("D26").Value = "'26.37"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "'0.163"
$ws.Range("E28").Value = "  +16.05%  "
$ws.Range("D29").Value = "'10.22"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("D31").Value = "'6.28"
$ws.Range("E31").Value = "  +6.79%  "
$ws.Range("D32").Value = "'34.87"
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("D33").Value = "'51.57"
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("D34").Value = "'0.0449"
$ws.Range("E34").Value = "  -9.32%  "
$ws.Range("D35").Value = "'0.0837"
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("D36").Value = "'5.12"
$ws.Range("E36").Value = "  -5.50%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "'18.68"
$ws.Range("E38").Value = "  +2.77%  "
$ws.Range("D39").Value = "'3.13"
$ws.Range("E39").Value = "  -2.96%  "
$ws.Range("D40").Value = "'1.95"
$ws.Range("E40").Value = "  -3.04%  "
$ws.Range("E41").Value = "  +3.22%  "
$ws.Range("D43").Value = "'2.22"
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("D44").Value = "'119.18"
$ws.Range("E44").Value = "  -6.05%  "
$ws.Range("D45").Value = "'21.53"
$ws.Range("E45").Value = "  -7.37%  "
$ws.Range("D46").Value = "2.083.89"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").Value = "'3.27"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("D49").Value = "'0.945"
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("D50").Value = "'5.58"
$ws.Range("E50").Value = "  -5.69%  "
$ws.Range("D51").Value = "'0.191"
$ws.Range("E51").Value = "  +1.12%  "
